# Begin modification to have MCSB logo:
# - Column F header "MCSB" -> "isMCSB"
# - Fill in F2:F70 with 0 (new boolean-ish flag column, not yet wired up)
# - Move the active selection to H66 (where the work left off)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new/renamed column.
$ws.Range("F1").Value = "isMCSB"

# Populate the flag column for every data row with a default of 0.
$lastRow = 70
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Leave the selection where editing stopped.
$ws.Range("H66").Select()
